# Append the latest day's COVID-19 data row to the "Tabela1" table on the
# "Covid-19 podatki" sheet, then move the selection the way the original
# author's session left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table currently spans A1:J89 (88 data rows). Add a new row to the
# table so it grows to A1:J90 and the new row inherits list-object wiring
# (auto filter, styling, dimension) the way typing under the table would.
$lo = $ws.ListObjects.Item(1)
$newListRow = $lo.ListRows.Add()

# Carry over the formatting used by the rows immediately above (row 85 is
# the last row using the regular/non-highlighted banding before the
# bordered rows 86-89) so the new row looks consistent with its neighbours.
$ws.Range("A85:J85").Copy() | Out-Null
$ws.Range("A90:J90").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New day's data (2020-06-08 -> serial date 43990).
$ws.Range("A90").Value = 43990
$ws.Range("B90").Value = 84130
$ws.Range("C90").Value = 814
$ws.Range("D90").Value = 1486
$ws.Range("E90").Value = 1
$ws.Range("F90").Value = 6
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 109
$ws.Range("J90").Value = 0

# Leave the selection where the author's editing session ended up.
$ws.Range("E94").Select() | Out-Null
